$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.135.19"
$ws.Range("E2").Value = "  -3.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.77"
$ws.Range("E3").Value = "  -4.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -1.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.87"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("E6").Value = "  -1.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4661"
$ws.Range("E7").Value = "  -6.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4004"
$ws.Range("E8").Value = "  -4.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.98"
$ws.Range("E9").Value = "  -4.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08404"
$ws.Range("E10").Value = "  -5.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.043"
$ws.Range("E11").Value = "  -4.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.08"
$ws.Range("E12").Value = "  -3.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.894.67"
$ws.Range("E13").Value = "  -5.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.419"
$ws.Range("E14").Value = "  -7.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.050"
$ws.Range("E15").Value = "  -5.67%  "

$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.46"
$ws.Range("E17").Value = "  -3.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001066"
$ws.Range("E18").Value = "  -3.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06634"
$ws.Range("E19").Value = "  -1.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.90"
$ws.Range("E20").Value = "  -7.96%  "

$ws.Range("E21").Value = "  -1.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.737"
$ws.Range("E22").Value = "  -4.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.137.82"
$ws.Range("E23").Value = "  -3.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.17"
$ws.Range("E24").Value = "  -6.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.300"
$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.148.82"
$ws.Range("E26").Value = "  -4.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.16"
$ws.Range("E27").Value = "  -2.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.01"
$ws.Range("E28").Value = "  -3.85%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.136"
$ws.Range("E29").Value = "  -4.92%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.753"
$ws.Range("E30").Value = "  -8.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.35"
$ws.Range("E31").Value = "  -2.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9749"
$ws.Range("E32").Value = "  -6.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09632"
$ws.Range("E33").Value = "  -2.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.437"
$ws.Range("E34").Value = "  -6.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.647"
$ws.Range("E35").Value = "  -2.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.541"
$ws.Range("E36").Value = "  -4.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.831"
$ws.Range("E37").Value = "  -2.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.264"
$ws.Range("E38").Value = "  -3.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02294"
$ws.Range("E39").Value = "  -5.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06164"
$ws.Range("E40").Value = "  -3.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6162"
$ws.Range("E41").Value = "  -4.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.01"
$ws.Range("E42").Value = "  -4.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  -1.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1900"
$ws.Range("E44").Value = "  -3.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.314"
$ws.Range("E45").Value = "  -3.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5850"
$ws.Range("E46").Value = "  -5.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.78"
$ws.Range("E47").Value = "  -3.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.020"
$ws.Range("E48").Value = "  -6.76%  "

$ws.Range("E49").Value = "  -1.55%  "

$ws.Range("E50").Value = "  -0.54%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.43"
$ws.Range("E51").Value = "  -1.25%  "

